# "Generate Report for Archive"
# The localization status moved from "Ready for handoff" to "In Translation"
# for both tracked files/languages. Update every cell that shows that status
# (the Overview roll-up sheet's per-language columns, plus the Status column
# on each language detail sheet), then let the now-shorter text's columns
# shrink to fit, matching the narrower column widths Excel recorded when the
# report was regenerated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn) and F (de-de) hold the status per language.
$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus
$overview.Range("E3").Value2 = $newStatus
$overview.Range("F3").Value2 = $newStatus

# Per-language detail sheets: column C is "Status".
$zhcn.Range("C2").Value2 = $newStatus
$zhcn.Range("C3").Value2 = $newStatus

$dede.Range("C2").Value2 = $newStatus
$dede.Range("C3").Value2 = $newStatus

# The status column(s) narrow now that the text is shorter than
# "Ready for handoff" - mirror Excel's recalculated column widths.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
